$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.3410236666666666
$ws.Range("H2").Value = 1.023071
$ws.Range("I2").Value = 0.01850325494520333
$ws.Range("J2").Value = 0.01850325494520333
$ws.Range("M2").Value = 0.5373756666666667
$ws.Range("N2").Value = 1.612127
$ws.Range("O2").Value = 0.007472820128982582
$ws.Range("P2").Value = 0.007472820128982581
$ws.Range("Q2").Value = 0.1832578202241111
$ws.Range("R2").Value = 1.649320382017
$ws.Range("S2").Value = 0.0001382714960062119
$ws.Range("T2").Value = 0.0001382714960062119
$ws.Range("G3").Value = 0.3410236666666666
$ws.Range("H3").Value = 1.023071
$ws.Range("I3").Value = 0.01850325494520333
$ws.Range("J3").Value = 0.01850325494520333
$ws.Range("O3").Value = 0.1537223653287423
$ws.Range("P3").Value = 0.1537223653287423
$ws.Range("Q3").Value = 3.769771666332778
$ws.Range("R3").Value = 33.927944996995
$ws.Range("S3").Value = 0.002844364116457404
$ws.Range("T3").Value = 0.002844364116457403
$ws.Range("G4").Value = 0.3410236666666666
$ws.Range("H4").Value = 1.023071
$ws.Range("I4").Value = 0.01850325494520333
$ws.Range("J4").Value = 0.01850325494520333
$ws.Range("M4").Value = 30.561198
$ws.Range("N4").Value = 91.683594
$ws.Range("O4").Value = 0.4249882340167162
$ws.Range("P4").Value = 0.4249882340167161
$ws.Range("Q4").Value = 10.422091799686
$ws.Range("R4").Value = 93.79882619717398
$ws.Range("S4").Value = 0.007863665642723032
$ws.Range("T4").Value = 0.007863665642723032
$ws.Range("G5").Value = 0.3410236666666666
$ws.Range("H5").Value = 1.023071
$ws.Range("I5").Value = 0.01850325494520333
$ws.Range("J5").Value = 0.01850325494520333
$ws.Range("M5").Value = 29.75783666666667
$ws.Range("N5").Value = 89.27351
$ws.Range("O5").Value = 0.4138165805255589
$ws.Range("P5").Value = 0.4138165805255589
$ws.Range("Q5").Value = 10.14812657213444
$ws.Range("R5").Value = 91.33313914920998
$ws.Range("S5").Value = 0.007656953690016679
$ws.Range("T5").Value = 0.007656953690016679
$ws.Range("I6").Value = 0.2085050756621187
$ws.Range("J6").Value = 0.2085050756621187
$ws.Range("M6").Value = 0.5373756666666667
$ws.Range("N6").Value = 1.612127
$ws.Range("O6").Value = 0.007472820128982582
$ws.Range("P6").Value = 0.007472820128982581
$ws.Range("Q6").Value = 2.065052110272556
$ws.Range("R6").Value = 18.585468992453
$ws.Range("S6").Value = 0.001558120926402917
$ws.Range("T6").Value = 0.001558120926402916
$ws.Range("I7").Value = 0.2085050756621187
$ws.Range("J7").Value = 0.2085050756621187
$ws.Range("O7").Value = 0.1537223653287423
$ws.Range("P7").Value = 0.1537223653287423
$ws.Range("S7").Value = 0.03205189341382927
$ws.Range("T7").Value = 0.03205189341382927
$ws.Range("I8").Value = 0.2085050756621187
$ws.Range("J8").Value = 0.2085050756621187
$ws.Range("M8").Value = 30.561198
$ws.Range("N8").Value = 91.683594
$ws.Range("O8").Value = 0.4249882340167162
$ws.Range("P8").Value = 0.4249882340167161
$ws.Range("Q8").Value = 117.441987676574
$ws.Range("R8").Value = 1056.977889089166
$ws.Range("S8").Value = 0.08861220388916562
$ws.Range("T8").Value = 0.0886122038891656
$ws.Range("I9").Value = 0.2085050756621187
$ws.Range("J9").Value = 0.2085050756621187
$ws.Range("M9").Value = 29.75783666666667
$ws.Range("N9").Value = 89.27351
$ws.Range("O9").Value = 0.4138165805255589
$ws.Range("P9").Value = 0.4138165805255589
$ws.Range("Q9").Value = 114.3547935224322
$ws.Range("R9").Value = 1029.19314170189
$ws.Range("S9").Value = 0.0862828574327209
$ws.Range("T9").Value = 0.08628285743272089
$ws.Range("G10").Value = 0.2092423333333333
$ws.Range("H10").Value = 0.627727
$ws.Range("I10").Value = 0.0113530661283407
$ws.Range("J10").Value = 0.0113530661283407
$ws.Range("M10").Value = 0.5373756666666667
$ws.Range("N10").Value = 1.612127
$ws.Range("O10").Value = 0.007472820128982582
$ws.Range("P10").Value = 0.007472820128982581
$ws.Range("Q10").Value = 0.1124417383698889
$ws.Range("R10").Value = 1.011975645329
$ws.Range("S10").Value = 0.00008483942108953475
$ws.Range("T10").Value = 0.00008483942108953473
$ws.Range("G11").Value = 0.2092423333333333
$ws.Range("H11").Value = 0.627727
$ws.Range("I11").Value = 0.0113530661283407
$ws.Range("J11").Value = 0.0113530661283407
$ws.Range("O11").Value = 0.1537223653287423
$ws.Range("P11").Value = 0.1537223653287423
$ws.Range("Q11").Value = 2.313023689257223
$ws.Range("R11").Value = 20.817213203315
$ws.Range("S11").Value = 0.00174522017898216
$ws.Range("T11").Value = 0.001745220178982159
$ws.Range("G12").Value = 0.2092423333333333
$ws.Range("H12").Value = 0.627727
$ws.Range("I12").Value = 0.0113530661283407
$ws.Range("J12").Value = 0.0113530661283407
$ws.Range("M12").Value = 30.561198
$ws.Range("N12").Value = 91.683594
$ws.Range("O12").Value = 0.4249882340167162
$ws.Range("P12").Value = 0.4249882340167161
$ws.Range("Q12").Value = 6.394696378982
$ws.Range("R12").Value = 57.55226741083801
$ws.Range("S12").Value = 0.004824919524558513
$ws.Range("T12").Value = 0.004824919524558512
$ws.Range("G13").Value = 0.2092423333333333
$ws.Range("H13").Value = 0.627727
$ws.Range("I13").Value = 0.0113530661283407
$ws.Range("J13").Value = 0.0113530661283407
$ws.Range("M13").Value = 29.75783666666667
$ws.Range("N13").Value = 89.27351
$ws.Range("O13").Value = 0.4138165805255589
$ws.Range("P13").Value = 0.4138165805255589
$ws.Range("Q13").Value = 6.226599179085555
$ws.Range("R13").Value = 56.03939261177
$ws.Range("S13").Value = 0.004698087003710496
$ws.Range("T13").Value = 0.004698087003710496
$ws.Range("G14").Value = 14.03735666666667
$ws.Range("H14").Value = 42.11207
$ws.Range("I14").Value = 0.7616386032643372
$ws.Range("J14").Value = 0.7616386032643372
$ws.Range("M14").Value = 0.5373756666666667
$ws.Range("N14").Value = 1.612127
$ws.Range("O14").Value = 0.007472820128982582
$ws.Range("P14").Value = 0.007472820128982581
$ws.Range("Q14").Value = 7.543333896987779
$ws.Range("R14").Value = 67.89000507289
$ws.Range("S14").Value = 0.005691588285483918
$ws.Range("T14").Value = 0.005691588285483917
$ws.Range("G15").Value = 14.03735666666667
$ws.Range("H15").Value = 42.11207
$ws.Range("I15").Value = 0.7616386032643372
$ws.Range("J15").Value = 0.7616386032643372
$ws.Range("O15").Value = 0.1537223653287423
$ws.Range("P15").Value = 0.1537223653287423
$ws.Range("Q15").Value = 155.1728944487945
$ws.Range("R15").Value = 1396.55605003915
$ws.Range("S15").Value = 0.1170808876194735
$ws.Range("T15").Value = 0.1170808876194735
$ws.Range("G16").Value = 14.03735666666667
$ws.Range("H16").Value = 42.11207
$ws.Range("I16").Value = 0.7616386032643372
$ws.Range("J16").Value = 0.7616386032643372
$ws.Range("M16").Value = 30.561198
$ws.Range("N16").Value = 91.683594
$ws.Range("O16").Value = 0.4249882340167162
$ws.Range("P16").Value = 0.4249882340167161
$ws.Range("Q16").Value = 428.9984364866201
$ws.Range("R16").Value = 3860.98592837958
$ws.Range("S16").Value = 0.323687444960269
$ws.Range("T16").Value = 0.323687444960269
$ws.Range("G17").Value = 14.03735666666667
$ws.Range("H17").Value = 42.11207
$ws.Range("I17").Value = 0.7616386032643372
$ws.Range("J17").Value = 0.7616386032643372
$ws.Range("M17").Value = 29.75783666666667
$ws.Range("N17").Value = 89.27351
$ws.Range("O17").Value = 0.4138165805255589
$ws.Range("P17").Value = 0.4138165805255589
$ws.Range("Q17").Value = 417.7213669184111
$ws.Range("R17").Value = 3759.4923022657
$ws.Range("S17").Value = 0.3151786823991108
$ws.Range("T17").Value = 0.3151786823991108
